$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# ---------------------------------------------------------------------------
# 1) Row 5: AC5 gains a new label (new shared string), copying the header
#    formatting used by its neighbouring label cells (e.g. Y5).
# ---------------------------------------------------------------------------
$ws.Range("Y5").Copy()
$ws.Range("AC5").PasteSpecial(-4122) | Out-Null
$ws.Range("AC5").Value2 = "Apresentação pré feira com slides"

# ---------------------------------------------------------------------------
# 2) Rows 6, 9, 11, 12, 13: the label that used to sit in column Z was moved
#    two columns to the right, into column AB (format + text travel
#    together; the old Z cell becomes completely blank).
# ---------------------------------------------------------------------------
$rowsWithLabelMove = 6, 9, 11, 12, 13
foreach ($r in $rowsWithLabelMove) {
    $srcCell = $ws.Range("Z$r")
    $dstCell = $ws.Range("AB$r")
    $labelText = $srcCell.Value2

    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122) | Out-Null
    $dstCell.Value2 = $labelText

    $srcCell.Clear()
}

# ---------------------------------------------------------------------------
# 3) Rows 7, 8, 10: column Z held no content, only formatting - it was
#    cleared out entirely (no more style override on those cells).
# ---------------------------------------------------------------------------
$rowsWithBlankClear = 7, 8, 10
foreach ($r in $rowsWithBlankClear) {
    $ws.Range("Z$r").Clear()
}

# ---------------------------------------------------------------------------
# 4) Update the view: scroll/select near the newly edited area.
# ---------------------------------------------------------------------------
$ws.Range("AC6").Select()
